# RTM.xlsx — build out the Requirements Traceability Matrix sheet:
# header row (CRS / CYRS / SRS or HIS) + three requirement rows, with
# column widths, row heights, yellow header fill, wrap text and a
# bordered table layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- cell values (order matters: this is also the shared-strings order) ----
$ws.Range("A1").Value = "CRS"
$ws.Range("B1").Value = "CYRS"
$ws.Range("C1").Value = "SRS or HIS"

$ws.Range("A2").Value = "At startup, WELCOME mode shall be one of the following modes:`n" + `
    "● First mode: LEDs shall be ON from L6 to L1, then from R1 to R6 and vice versa, and then all LEDs are ON and OFF.`n" + `
    "● Second mode: LEDS from R1 to R6 are ON LED by LED and also the left branch at the same time, and then repeat the scenario again.`n"

$ws.Range("A3").Value = "Tail function shall be activated according to Tail switch."

$ws.Range("A4").Value = "TI function shall be activated be activated according to TI switch, LEDs shall be activated LED by LED from R1 to R6 or from L1 to L6."

# ---- column widths ----
$ws.Columns.Item(1).ColumnWidth = 116.140625
$ws.Columns.Item(2).ColumnWidth = 11.42578125
$ws.Columns.Item(3).ColumnWidth = 13.85546875

# ---- row heights ----
$ws.Rows.Item(1).RowHeight = 15.75
$ws.Rows.Item(2).RowHeight = 115.5
$ws.Rows.Item(3).RowHeight = 29.25
$ws.Rows.Item(4).RowHeight = 33

# ---- header row fill (yellow) ----
$ws.Range("A1:C1").Interior.Color = 65535

# ---- wrap text on the big requirement cell ----
$ws.Range("A2").WrapText = $true

# ---- borders: thin grid inside the header row ----
$ws.Range("A1:C1").Borders.Item(11).Weight = 2

# ---- borders: medium line separating header row from the data rows ----
$ws.Range("A1:C2").Borders.Item(12).Weight = -4138

# ---- borders: medium outer box around the whole table ----
$ws.Range("A1:C4").Borders.Item(7).Weight = -4138
$ws.Range("A1:C4").Borders.Item(10).Weight = -4138
$ws.Range("A1:C4").Borders.Item(8).Weight = -4138
$ws.Range("A1:C4").Borders.Item(9).Weight = -4138

# ---- selection, matching the saved cursor position ----
[void]$ws.Range("A2").Select()
